# Auto-generated edit script: updates cryptos list price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.944.75"
$ws.Range("D2").Style = $s
$ws.Range("E2").Value = "  -1.16%  "

$s = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.817.98"
$ws.Range("D3").Style = $s
$ws.Range("E3").Value = "  -0.06%  "

$ws.Range("E4").Value = "  -0.05%  "

$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.04"
$ws.Range("D5").Style = $s
$ws.Range("E5").Value = "  -1.12%  "

$ws.Range("E6").Value = "  -0.09%  "

$s = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4642"
$ws.Range("D7").Style = $s
$ws.Range("E7").Value = "  -0.68%  "

$s = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3707"
$ws.Range("D8").Style = $s
$ws.Range("E8").Value = "  -1.76%  "

$s = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07361"
$ws.Range("D9").Style = $s
$ws.Range("E9").Value = "  -0.62%  "

$s = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8731"
$ws.Range("D10").Style = $s
$ws.Range("E10").Value = "  +0.20%  "

$s = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.48"
$ws.Range("D11").Style = $s
$ws.Range("E11").Value = "  -0.54%  "

$s = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.858.23"
$ws.Range("D12").Style = $s
$ws.Range("E12").Value = "  +2.10%  "

$ws.Range("E13").Value = "  -1.03%  "

$s = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.522"
$ws.Range("D14").Style = $s
$ws.Range("E14").Value = "  -2.43%  "

$s = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07075"
$ws.Range("D15").Style = $s
$ws.Range("E15").Value = "  -0.10%  "

$s = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.68"
$ws.Range("D16").Style = $s
$ws.Range("E16").Value = "  -0.66%  "

$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("E20").Value = "  -1.33%  "

$s = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.954.51"
$ws.Range("D21").Style = $s
$ws.Range("E21").Value = "  -1.18%  "

$s = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.323"
$ws.Range("D22").Style = $s
$ws.Range("E22").Value = "  +0.15%  "

$s = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.60"
$ws.Range("D23").Style = $s
$ws.Range("E23").Value = "  -3.21%  "

$s = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.065.82"
$ws.Range("D24").Style = $s
$ws.Range("E24").Value = "  +0.80%  "

$s = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.908"
$ws.Range("D25").Style = $s
$ws.Range("E25").Value = "  -1.71%  "

$s = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.27"
$ws.Range("D26").Style = $s
$ws.Range("E26").Value = "  +0.44%  "

$s = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.44"
$ws.Range("D27").Style = $s
$ws.Range("E27").Value = "  -0.64%  "

$s = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.153"
$ws.Range("D28").Style = $s
$ws.Range("E28").Value = "  -4.28%  "

$s = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.318"
$ws.Range("D29").Style = $s
$ws.Range("E29").Value = "  +0.14%  "

$s = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.83"
$ws.Range("D30").Style = $s
$ws.Range("E30").Value = "  -1.21%  "

$s = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08920"
$ws.Range("D31").Style = $s
$ws.Range("E31").Value = "  -0.26%  "

$s = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7590"
$ws.Range("D32").Style = $s
$ws.Range("E32").Value = "  -3.10%  "

$s = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.157"
$ws.Range("D33").Style = $s
$ws.Range("E33").Value = "  -1.93%  "

$s = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.477"
$ws.Range("D34").Style = $s
$ws.Range("E34").Value = "  -1.00%  "

$s = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.923"
$ws.Range("D35").Style = $s
$ws.Range("E35").Value = "  -0.42%  "

$ws.Range("E36").Value = "  -0.04%  "

$s = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.094"
$ws.Range("D37").Style = $s
$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("E38").Value = "  -0.43%  "

$s = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05262"
$ws.Range("D39").Style = $s
$ws.Range("E39").Value = "  +0.26%  "

$s = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.938"
$ws.Range("D40").Style = $s
$ws.Range("E40").Value = "  +1.84%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$s = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.264"
$ws.Range("D41").Style = $s
$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$s = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.391"
$ws.Range("D42").Style = $s
$ws.Range("E42").Value = "  +1.50%  "

$s = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5344"
$ws.Range("D43").Style = $s
$ws.Range("E43").Value = "  +0.38%  "

$s = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1669"
$ws.Range("D44").Style = $s

$ws.Range("E45").Value = "  -1.62%  "

$s = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4961"
$ws.Range("D46").Style = $s
$ws.Range("E46").Value = "  -2.11%  "

$s = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.38"
$ws.Range("D47").Style = $s
$ws.Range("E47").Value = "  -0.57%  "

$s = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.679"
$ws.Range("D48").Style = $s
$ws.Range("E48").Value = "  +0.63%  "

$s = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("D49").Style = $s
$ws.Range("E49").Value = "  -0.03%  "

$s = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.34"
$ws.Range("D50").Style = $s
$ws.Range("E50").Value = "  -2.08%  "

$s = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06290"
$ws.Range("D51").Style = $s
$ws.Range("E51").Value = "  -0.71%  "
